$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36 - this pushes the existing rows
# 36..87 down to 37..88 (matching Excel's native Rows.Insert behaviour,
# including carrying the date number format on column D down to the new
# row and bumping the sheet's used-range dimension automatically).
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly record.
$ws.Cells.Item(36, 1).Value  = 10
$ws.Cells.Item(36, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(36, 3).Value  = 'La Araucanía'
$ws.Cells.Item(36, 4).Value  = 44775
$ws.Cells.Item(36, 5).Value  = 9
$ws.Cells.Item(36, 6).Value  = 'Fruta'
$ws.Cells.Item(36, 7).Value  = 100108
$ws.Cells.Item(36, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(36, 9).Value  = 100108007
$ws.Cells.Item(36, 10).Value = 'Coco'
$ws.Cells.Item(36, 11).Value = 'Sin especificar'
$ws.Cells.Item(36, 12).Value = 'Primera'
$ws.Cells.Item(36, 13).Value = 20
$ws.Cells.Item(36, 14).Value = 30000
$ws.Cells.Item(36, 15).Value = 30000
$ws.Cells.Item(36, 16).Value = 30000
$ws.Cells.Item(36, 17).Value = '$/malla 20 unidades'
$ws.Cells.Item(36, 18).Value = 'Perú'
$ws.Cells.Item(36, 19).Value = 1500
$ws.Cells.Item(36, 20).Value = 20
